$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet in the right spot.
#    We duplicate the existing "总计" sheet (so it keeps the same sheetPr /
#    sheetView / pageMargins structure), rename the *original* worksheet
#    object to "2022-Q1" (it keeps its original sheetId/part) and rename the
#    freshly created copy back to "总计" (it becomes the new part for the
#    totals sheet). This reproduces the exact sheetId / part assignment of
#    the target workbook.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Copy($null, $total)
$total.Name = "2022-Q1"
$newTotal = $wb.Worksheets.Item("总计 (2)")
$newTotal.Name = "总计"

# ---------------------------------------------------------------------------
# 2) Rewrite the "2022-Q1" sheet with the new fund-holdings table.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-Q1")

# Extend the header styling (bold + border, same as the existing B1:D1
# cells) onto the three brand-new header cells E1:H1.
$ws.Range("D1").Copy($ws.Range("E1:H1"))

# Extend the row-index styling (bold + border, same as the existing A2
# cell) onto the brand-new rows A6:A14.
$ws.Range("A2").Copy($ws.Range("A6:A14"))

# Header row.
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

# Columns B, D, E, F all hold numeric-looking text that must stay text
# (t="inlineStr"/shared-string), so force the cells to Text format before
# assigning the values (otherwise Excel auto-converts them to numbers).
$ws.Range("B2:B14").NumberFormat = "@"
$ws.Range("D2:G14").NumberFormat = "@"

$rows = @(
    @(0, "002601", "中银证券价值精选灵活配置混合", "3.41", "93.74", "4.89", "0.1667", 7),
    @(1, "561550", "华泰柏瑞中证500增强策略ETF", "9.68", "98.93", "1.12", "0.1084", 6),
    @(2, "011269", "中银证券优势制造股票型证券投资基金A", "1.39", "93.51", "6.55", "0.0910", 2),
    @(3, "003655", "信达澳银新财富灵活配置混合", "11.86", "25.86", "0.65", "0.0771", 8),
    @(4, "008997", "同泰竞争优势混合A", "1.33", "94.35", "3.08", "0.0410", 10),
    @(5, "008998", "同泰竞争优势混合C", "0.95", "94.35", "3.08", "0.0293", 10),
    @(6, "011934", "中航量化阿尔法六个月持有股票A", "2.96", "90.35", "0.96", "0.0284", 3),
    @(7, "011935", "中航量化阿尔法六个月持有股票C", "2.19", "90.35", "0.96", "0.0210", 3),
    @(8, "011270", "中银证券优势制造股票型证券投资基金C", "0.21", "93.51", "6.55", "0.0138", 2),
    @(9, "004192", "招商中证500指数增强A", "0.96", "94.32", "0.90", "0.0086", 10),
    @(10, "004193", "招商中证500指数增强C", "0.42", "94.32", "0.90", "0.0038", 10),
    @(11, "003586", "先锋精一灵活配置混合A", "0.03", "92.66", "4.35", "0.0013", 5),
    @(12, "003587", "先锋精一灵活配置混合C", "0.03", "92.66", "4.35", "0.0013", 5)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Update the "总计" sheet: push the existing rows down by one and insert
#    the new 2022-Q1 summary row at the top of the data (row 2).
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Extend the row-index styling (bold + border, same as the existing A2:A5
# cells) onto the brand-new row A6 that appears once everything shifts down.
$wsTotal.Range("A5").Copy($wsTotal.Range("A6"))

for ($row = 5; $row -ge 2; $row--) {
    $dst = $row + 1
    $wsTotal.Cells.Item($dst, 1).Value = $row - 1
    $wsTotal.Cells.Item($dst, 2).Value = $wsTotal.Cells.Item($row, 2).Text
    $wsTotal.Cells.Item($dst, 3).Value = $wsTotal.Cells.Item($row, 3).Text
    $wsTotal.Cells.Item($dst, 4).Value = $wsTotal.Cells.Item($row, 4).Text
}

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 13
$wsTotal.Cells.Item(2, 4).Value = 0.59
